# Adds Cat's phase/native-status codes ("I" = Introduced, "N" = Native,
# "NA" = Not Applicable) to column F for the rows that were previously
# missing a value in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    7  = "I"
    19 = "N"
    35 = "N"
    38 = "N"
    50 = "N"
    55 = "N"
    56 = "I"
    79 = "I"
    83 = "N"
    84 = "NA"
    85 = "NA"
    86 = "NA"
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

# Reflect the saved view state: selection moved to F1 (and the sheet
# scrolled back to the top, since the file no longer pins topLeftCell).
$ws.Range("F1").Select()
